$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 45; $row++) {
    $qCell = $ws.Cells.Item($row, 17)   # Column Q (season_x)
    $sCell = $ws.Cells.Item($row, 19)   # Column S (season_ending_year_y)

    $qCell.Value2 = $qCell.Value2 - 1
    $sCell.Value2 = $sCell.Value2 + 1
}
